$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.834999999999999
$ws.Range("B4").Value = 6.578
$ws.Range("E6").Value = 12.384
$ws.Range("B7").Value = 7.478
$ws.Range("E7").Value = 12.405
$ws.Range("B8").Value = 7.355
$ws.Range("E8").Value = 12.087
$ws.Range("A11").Value = -21.556
$ws.Range("A12").Value = -21.36
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 6.531999999999999
$ws.Range("A15").Value = -21.221
$ws.Range("E19").Value = 12.246
$ws.Range("E21").Value = 13.144
$ws.Range("B22").Value = 6.928999999999999
$ws.Range("E24").Value = 12.57
$ws.Range("E25").Value = 12.246
